$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: phone number A43 was stored as text; it should become a real number.
$ws.Range("A43").Value = 79174445

# Row 44: new payment record 79174445 (Cash) 2025-08-18T17:42:58
# A44 keeps the phone number as text (matches how it originally looked in A43
# before the fix above), so force a text format before assigning it.
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "79174445"
$ws.Range("A44").Style = "Normal"

$ws.Range("B44").Value = ""
$ws.Range("C44").Value = "Cash"
$ws.Range("D44").Value = "2025-08-18T17:42:58"
$ws.Range("E44").Value = 60
$ws.Range("F44").Value = ""
$ws.Range("G44").Value = 60
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
